# The order-item catalog must no longer allow deleting items/customers once
# an order exists for them. The stray duplicate/test row "aaa" (id 1159) is
# removed from the Articulos sheet, shifting every following row up by one,
# and the former placeholder row "jmaon" (now the last row) is renamed to
# "submit".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire row for article "aaa" (row 205) - this shifts all rows
# below it up by one, matching the diff exactly.
$ws.Rows.Item(205).Delete()

# The former last row ("jmaon") is now row 212 after the shift above;
# rename its article name to "submit" while keeping the rest of the row.
$ws.Range("B212").Value = "submit"
